$wb = $excel.ActiveWorkbook

# Sheet "VENTAS POR GRUPO": PORCELANATO total for HIDALGO HIDALGO PEDRO GUSTAVO (row 14)
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M14").Value = 5583.28

# Sheet "VENTA MENSUAL": septiembre value for the same advisor (row 14) and the column total (row 23)
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F14").Value = 5583.28
$wsMensual.Range("F23").Value = 25671.7

# Sheet "CUMPLIMIENTO MENSUAL": PORCELANATO row (12) and TOTAL row (15) for the advisor
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D12").Value = 23578.45
$wsCumplimiento.Range("E12").Value = 13245.1930921171
$wsCumplimiento.Range("F12").Value = 0.64030736831271

$wsCumplimiento.Range("D15").Value = 25671.7
$wsCumplimiento.Range("E15").Value = 29753.04316613378
$wsCumplimiento.Range("F15").Value = 0.4631812171515158
